# Add additional signal documentation rows to the "Signals" worksheet.
# Cell values are written in the specific sequence below so that newly
# introduced shared-string entries are appended to xl/sharedStrings.xml
# in the same order as the target revision.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Signals")

# Row 22 - CYCLE_START
$ws.Range("C22").Value = "CYCLE_START"
$ws.Range("D22").Value = "CNC Cycle Start"
$ws.Range("E22").Value = "DC"
$ws.Range("F22").Value = "24V"

# Row 23 - FEED_OVER (tag + voltage type first)
$ws.Range("C23").Value = "FEED_OVER"
$ws.Range("E23").Value = "Analog"

# Row 24 - RAPID_OVER
$ws.Range("C24").Value = "RAPID_OVER"
$ws.Range("D24").Value = "Rapid Speed Override"
$ws.Range("E24").Value = "Analog"

# Back to row 23 to fill in the description
$ws.Range("D23").Value = "Feed Speed Override"

# Row 25 - SPINDLE_ENC
$ws.Range("C25").Value = "SPINDLE_ENC"
$ws.Range("D25").Value = "Spindle Encoder"
$ws.Range("E25").Value = "MUX"

# Drive alarms - tags entered Z, X, Y ...
$ws.Range("C28").Value = "ALARM_Z"
$ws.Range("C26").Value = "ALARM_X"
$ws.Range("C27").Value = "ALARM_Y"

# ... then descriptions entered X, Y, Z
$ws.Range("D26").Value = "Drive Alarm X-Axis"
$ws.Range("D27").Value = "Drive Alarm Y-Axis"
$ws.Range("D28").Value = "Drive Alarm Z-Axis"

# Voltage columns for the alarm rows
$ws.Range("E26").Value = "DC"
$ws.Range("F26").Value = "24V"
$ws.Range("E27").Value = "DC"
$ws.Range("F27").Value = "24V"
$ws.Range("E28").Value = "DC"
$ws.Range("F28").Value = "24V"

# Row 29 - PENDANT
$ws.Range("C29").Value = "PENDANT"
$ws.Range("D29").Value = "Motion Pendant"
$ws.Range("E29").Value = "MUX"

$ws.Range("C30").Select() | Out-Null
